$d = $word.ActiveDocument

$replacements = @(
    @("45×36=1620", "76×86=6536"),
    @("19×56=1064", "95×66=6270"),
    @("88×65=5720", "53×69=3657"),
    @("74×83=6142", "67×90=6030"),
    @("49×99=4851", "35×70=2450"),
    @("95×86=8170", "54×67=3618"),
    @("23×81=1863", "53×85=4505"),
    @("49×29=1421", "95×29=2755"),
    @("37×11=407",  "26×59=1534"),
    @("94×56=5264", "79×90=7110"),
    @("99×77=7623", "38×38=1444"),
    @("89×41=3649", "36×83=2988"),
    @("19×62=1178", "20×47=940"),
    @("42×17=714",  "14×97=1358"),
    @("88×84=7392", "86×62=5332"),
    @("55×74=4070", "82×23=1886"),
    @("12×93=1116", "97×68=6596"),
    @("72×60=4320", "90×97=8730"),
    @("27×28=756",  "83×81=6723"),
    @("61×55=3355", "59×50=2950"),
    @("74×36=2664", "66×86=5676"),
    @("36×97=3492", "97×48=4656"),
    @("30×11=330",  "78×44=3432"),
    @("44×20=880",  "93×17=1581"),
    @("72×55=3960", "86×12=1032")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
